# Auto-generated edits: refresh market-price derived columns (H-N)
# across all 8 job sheets, matching the scheduled-runner price update.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 745.9091
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 770.5
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 770.5
$ws.Range("M19").Value = -325
$ws.Range("N19").Value = -1120.5
$ws.Range("H31").Value = 728.5
$ws.Range("I31").Value = 728.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2185.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1955.5
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 71430776
$ws.Range("I32").Value = 1599.5
$ws.Range("J32").Value = 100002450
$ws.Range("K32").Value = 1599.5
$ws.Range("L32").Value = 100002450
$ws.Range("M32").Value = -1273.5
$ws.Range("N32").Value = -100003102
$ws.Range("H33").Value = 15637.41
$ws.Range("I33").Value = 18183.281
$ws.Range("K33").Value = 18183.281
$ws.Range("M33").Value = -17954.281
$ws.Range("H64").Value = 6136.6
$ws.Range("I64").Value = 6136.6
$ws.Range("K64").Value = 6136.6
$ws.Range("M64").Value = -5888.6
$ws.Range("H67").Value = 6136.6
$ws.Range("I67").Value = 6136.6
$ws.Range("K67").Value = 6136.6
$ws.Range("M67").Value = -5278.6
$ws.Range("H69").Value = 17945.1
$ws.Range("I69").Value = 11653
$ws.Range("J69").Value = 20641.715
$ws.Range("K69").Value = 34959
$ws.Range("L69").Value = 61925.145
$ws.Range("M69").Value = -34085
$ws.Range("N69").Value = -63673.145
$ws.Range("H70").Value = 1377
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 1546.25
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 4638.75
$ws.Range("M70").Value = -1830
$ws.Range("N70").Value = -5178.75
$ws.Range("H72").Value = 17945.1
$ws.Range("I72").Value = 11653
$ws.Range("J72").Value = 20641.715
$ws.Range("K72").Value = 104877
$ws.Range("L72").Value = 185775.435
$ws.Range("M72").Value = -100509
$ws.Range("N72").Value = -194511.435
$ws.Range("H73").Value = 1377
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 1546.25
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 4638.75
$ws.Range("M73").Value = -1164
$ws.Range("N73").Value = -6510.75
$ws.Range("H103").Value = 529.6923
$ws.Range("I103").Value = 384.2857
$ws.Range("J103").Value = 699.3333
$ws.Range("K103").Value = 1152.8571
$ws.Range("L103").Value = 2097.9999
$ws.Range("M103").Value = -566.8571000000002
$ws.Range("N103").Value = -3269.9999
$ws.Range("H116").Value = 4607.5
$ws.Range("I116").Value = 4715.909
$ws.Range("K116").Value = 4715.909
$ws.Range("M116").Value = -1273.909
$ws.Range("H125").Value = 1678.2858
$ws.Range("I125").Value = 736.8333
$ws.Range("K125").Value = 6631.4997
$ws.Range("M125").Value = -4171.4997
$ws.Range("H137").Value = 9092191
$ws.Range("I137").Value = 1261.875
$ws.Range("K137").Value = 3785.625
$ws.Range("M137").Value = -1235.625

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1693.7059
$ws.Range("I2").Value = 746.2
$ws.Range("K2").Value = 746.2
$ws.Range("M2").Value = -633.2
$ws.Range("H45").Value = 55994.5
$ws.Range("I45").Value = 114876.664
$ws.Range("J45").Value = 7818.1816
$ws.Range("K45").Value = 114876.664
$ws.Range("L45").Value = 7818.1816
$ws.Range("M45").Value = -114499.664
$ws.Range("N45").Value = -8572.1816
$ws.Range("H61").Value = 2256726.5
$ws.Range("I61").Value = 47864.652
$ws.Range("K61").Value = 47864.652
$ws.Range("M61").Value = -47652.652
$ws.Range("H97").Value = 6281.5293
$ws.Range("I97").Value = 6281.5293
$ws.Range("K97").Value = 6281.5293
$ws.Range("M97").Value = -5785.5293
$ws.Range("H110").Value = 1605.9474
$ws.Range("I110").Value = 1442.0588
$ws.Range("K110").Value = 1442.0588
$ws.Range("M110").Value = 602.9412
$ws.Range("H116").Value = 1693.7059
$ws.Range("I116").Value = 746.2
$ws.Range("K116").Value = 746.2
$ws.Range("M116").Value = 1547.8
$ws.Range("H122").Value = 1113.2858
$ws.Range("I122").Value = 1113.2858
$ws.Range("K122").Value = 3339.8574
$ws.Range("M122").Value = -889.8574000000003
$ws.Range("H132").Value = 2314.6155
$ws.Range("I132").Value = 2201.4
$ws.Range("J132").Value = 2692
$ws.Range("K132").Value = 6604.200000000001
$ws.Range("L132").Value = 8076
$ws.Range("M132").Value = -4074.200000000001
$ws.Range("N132").Value = -13136
$ws.Range("H136").Value = 2256726.5
$ws.Range("I136").Value = 47864.652
$ws.Range("K136").Value = 143593.956
$ws.Range("M136").Value = -141043.956

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1693.7059
$ws.Range("I3").Value = 746.2
$ws.Range("K3").Value = 746.2
$ws.Range("M3").Value = -632.2
$ws.Range("H86").Value = 4383
$ws.Range("I86").Value = 1999.6666
$ws.Range("K86").Value = 1999.6666
$ws.Range("M86").Value = -876.6666
$ws.Range("H89").Value = 4383
$ws.Range("I89").Value = 1999.6666
$ws.Range("K89").Value = 9998.333000000001
$ws.Range("M89").Value = -4382.333000000001
$ws.Range("H94").Value = 3227.9333
$ws.Range("I94").Value = 2709.2856
$ws.Range("J94").Value = 3681.75
$ws.Range("K94").Value = 2709.2856
$ws.Range("L94").Value = 3681.75
$ws.Range("M94").Value = -2258.2856
$ws.Range("N94").Value = -4583.75
$ws.Range("H99").Value = 7337.857
$ws.Range("I99").Value = 8451.529
$ws.Range("K99").Value = 8451.529
$ws.Range("M99").Value = -6953.529
$ws.Range("H107").Value = 6222.8604
$ws.Range("I107").Value = 7289.9116
$ws.Range("K107").Value = 7289.9116
$ws.Range("M107").Value = -5369.9116
$ws.Range("H134").Value = 18751956
$ws.Range("I134").Value = 1818.8334
$ws.Range("J134").Value = 75002370
$ws.Range("K134").Value = 5456.5002
$ws.Range("L134").Value = 225007110
$ws.Range("M134").Value = -2921.5002
$ws.Range("N134").Value = -225012180

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2328.6155
$ws.Range("I58").Value = 2555.8
$ws.Range("J58").Value = 2186.625
$ws.Range("K58").Value = 2555.8
$ws.Range("L58").Value = 2186.625
$ws.Range("M58").Value = -2352.8
$ws.Range("N58").Value = -2592.625
$ws.Range("H132").Value = 61615
$ws.Range("I132").Value = 69163.664
$ws.Range("K132").Value = 207490.992
$ws.Range("M132").Value = -204960.992
$ws.Range("H136").Value = 2328.6155
$ws.Range("I136").Value = 2555.8
$ws.Range("J136").Value = 2186.625
$ws.Range("K136").Value = 7667.400000000001
$ws.Range("L136").Value = 6559.875
$ws.Range("M136").Value = -5117.400000000001
$ws.Range("N136").Value = -11659.875

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 21482.416
$ws.Range("J41").Value = 41914.832
$ws.Range("L41").Value = 125744.496
$ws.Range("N41").Value = -126420.496

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H102").Value = 12821796
$ws.Range("I102").Value = 15626240
$ws.Range("K102").Value = 15626240
$ws.Range("M102").Value = -15624618
$ws.Range("H132").Value = 856590.9399999999
$ws.Range("I132").Value = 2426.2
$ws.Range("J132").Value = 1027423.9
$ws.Range("K132").Value = 7278.599999999999
$ws.Range("L132").Value = 3082271.7
$ws.Range("M132").Value = -4748.599999999999
$ws.Range("N132").Value = -3087331.7

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4170.2
$ws.Range("I22").Value = 1443.1
$ws.Range("J22").Value = 5988.2666
$ws.Range("K22").Value = 1443.1
$ws.Range("L22").Value = 5988.2666
$ws.Range("M22").Value = -1148.1
$ws.Range("N22").Value = -6578.2666
$ws.Range("H27").Value = 4170.2
$ws.Range("I27").Value = 1443.1
$ws.Range("J27").Value = 5988.2666
$ws.Range("K27").Value = 1443.1
$ws.Range("L27").Value = 5988.2666
$ws.Range("M27").Value = -1336.1
$ws.Range("N27").Value = -6202.2666
$ws.Range("H68").Value = 6273.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 6273.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H123").Value = 127495
$ws.Range("J123").Value = 127495
$ws.Range("L123").Value = 127495
$ws.Range("N123").Value = -137295

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3000000
$ws.Range("J54").Value = 3000000
$ws.Range("L54").Value = 3000000
$ws.Range("N54").Value = -3001040
$ws.Range("H96").Value = 22393.9
$ws.Range("I96").Value = 3378.8
$ws.Range("J96").Value = 41409
$ws.Range("K96").Value = 3378.8
$ws.Range("L96").Value = 41409
$ws.Range("M96").Value = -2005.8
$ws.Range("N96").Value = -44155
